{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst titleText =\n  \"Answers: Rearranging equations involving trigonometry and logarithms\";\nconst authorText = \"Ellie Gurini\";\nconst abstractText =\n  \"This is an answer set relating to the questions based on Guide, \" +\n  \"Introduction to rearranging equations involving trigonometry and logarithms.\";\n\nfor (const paragraph of paragraphs.items) {\n  const currentText = paragraph.text;\n  if (currentText === titleText && paragraph.style === \"Title\") {\n    paragraph.insertText(titleText, \"Replace\");\n  } else if (currentText === authorText && paragraph.style === \"Author\") {\n    paragraph.insertText(authorText, \"Replace\");\n  } else if (currentText === abstractText && paragraph.style === \"Abstract\") {\n    paragraph.insertText(abstractText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-DocText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Title paragraph: merge the word-by-word runs into a single run.\nReplace-DocText \"Answers: Rearranging equations involving trigonometry and logarithms\" \"Answers: Rearranging equations involving trigonometry and logarithms\"\n\n# Author paragraph: merge \"Ellie\" / \" \" / \"Gurini\" runs into a single run.\nReplace-DocText \"Ellie Gurini\" \"Ellie Gurini\"\n\n# Abstract paragraph: merge the word-by-word runs into a single run.\nReplace-DocText \"This is an answer set relating to the questions based on Guide, Introduction to rearranging equations involving trigonometry and logarithms.\" \"This is an answer set relating to the questions based on Guide, Introduction to rearranging equations involving trigonometry and logarithms.\"\n"}
